$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegDetails")

# Update FirstName / MiddleName/Initial / LastName values
$ws.Range("A2").Value = "Sam"
$ws.Range("B2").Value = "Ron"
$ws.Range("C2").Value = "Jin"

$ws.Range("A3").Value = "Peter"
$ws.Range("B3").Value = "Kio"
$ws.Range("C3").Value = "Conery"

$ws.Range("A4").Value = "Butna"
$ws.Range("B4").Value = "Amy"
$ws.Range("C4").Value = "Swan"

$ws.Range("A5").Value = "Jen"
$ws.Range("B5").Value = "Loper"
$ws.Range("C5").Value = "Kou"

# Update the displayed EmailAddress text (hyperlink target stays the same)
$ws.Range("D2").Value = "a120@email.com"
$ws.Range("D3").Value = "a121@email.com"
$ws.Range("D4").Value = "a122@email.com"
$ws.Range("D5").Value = "a123@email.com"

# Clear the Status column data cells (F2:F5) - header stays
$ws.Range("F2:F5").Clear()

# Update sheet view: scroll position (top-left visible column -> D) + selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("E8").Select()

# Update workbook window size
$excel.ActiveWindow.Width = 11505
$excel.ActiveWindow.Height = 4110
